# Updates the "想去人数" (interested-count) column F on the three sheets
# that carry data rows, re-syncing cached scrape numbers per commit
# "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value  = 3357
$ws.Range("F5").Value  = 217
$ws.Range("F6").Value  = 4850
$ws.Range("F7").Value  = 471
$ws.Range("F8").Value  = 301
$ws.Range("F9").Value  = 177
$ws.Range("F10").Value = 631
$ws.Range("F12").Value = 37
$ws.Range("F14").Value = 662
$ws.Range("F17").Value = 90
$ws.Range("F18").Value = 147
$ws.Range("F19").Value = 345
$ws.Range("F20").Value = 4767
$ws.Range("F21").Value = 24
$ws.Range("F24").Value = 5906
$ws.Range("F26").Value = 1200
$ws.Range("F27").Value = 246
$ws.Range("F28").Value = 674
$ws.Range("F29").Value = 4423
$ws.Range("F31").Value = 95
$ws.Range("F32").Value = 127
$ws.Range("F33").Value = 871
$ws.Range("F34").Value = 75
$ws.Range("F36").Value = 793
$ws.Range("F37").Value = 854
$ws.Range("F38").Value = 2

# Sheet "本地生活" (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 220
$ws.Range("F3").Value = 1092

# Sheet "全部类型" (All types) -- merged view of the other sheets
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 220
$ws.Range("F4").Value  = 1092
$ws.Range("F8").Value  = 3357
$ws.Range("F9").Value  = 217
$ws.Range("F10").Value = 4850
$ws.Range("F11").Value = 471
$ws.Range("F12").Value = 301
$ws.Range("F13").Value = 177
$ws.Range("F14").Value = 631
$ws.Range("F16").Value = 37
$ws.Range("F18").Value = 662
$ws.Range("F22").Value = 90
$ws.Range("F23").Value = 147
$ws.Range("F24").Value = 345
$ws.Range("F25").Value = 4767
$ws.Range("F26").Value = 24
$ws.Range("F29").Value = 5906
$ws.Range("F31").Value = 1200
$ws.Range("F32").Value = 246
$ws.Range("F33").Value = 674
$ws.Range("F34").Value = 4423
$ws.Range("F37").Value = 95
$ws.Range("F38").Value = 127
$ws.Range("F39").Value = 871
$ws.Range("F40").Value = 75
$ws.Range("F42").Value = 793
$ws.Range("F43").Value = 854
$ws.Range("F45").Value = 2
